# Apply updated cryptocurrency price/volume figures to Sheet1
# (matches the Sat Mar 23 23:56:38 UTC 2024 GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.351.16"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "'3.355.35"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'554.79"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Value = "'173.44"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").Value = "'3.348.60"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.626"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "'0.162"
$ws.Range("E11").Value = "  +6.60%  "
$ws.Range("D12").Value = "'53.53"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "'9.01"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "'3.894.57"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'18.12"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "'3.351.04"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'64.267.52"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "'11.67"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "'0.983"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'457.18"
$ws.Range("E22").Value = "  +7.50%  "
$ws.Range("D23").Value = "'4.84"
$ws.Range("E23").Value = "  +9.20%  "
$ws.Range("D25").Value = "'85.60"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "'13.49"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "'2.95"
$ws.Range("E27").Value = "  +8.12%  "
$ws.Range("D28").Value = "'10.70"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'8.66"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'30.33"
$ws.Range("E30").Value = "  +3.87%  "
$ws.Range("D31").Value = "'6.64"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").Value = "'11.37"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'570.45"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'60.89"
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'3.62"
$ws.Range("E37").Value = "  +5.06%  "
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("D39").Value = "'35.14"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.0₃0735"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'0.366"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'3.065.14"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'139.20"
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("D51").Value = "'8.10"
$ws.Range("E51").Value = "  +1.06%  "
